$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (unchanged, rewritten defensively) ---
$ws.Cells.Item(1,1).Value = "Conta"
$ws.Cells.Item(1,2).Value = "Nome"
$ws.Cells.Item(1,3).Value = "Saldo"

# --- Data rows (account numbers written with a leading apostrophe so the
#     9-digit, zero-padded account codes stay text instead of becoming numbers) ---
$ws.Cells.Item(2,1).Value = "'005995120"
$ws.Cells.Item(2,2).Value = "Erik"
$ws.Cells.Item(2,3).Value = 73064.56
$ws.Cells.Item(3,1).Value = "'004212476"
$ws.Cells.Item(3,2).Value = "Maria"
$ws.Cells.Item(3,3).Value = 59164.35
$ws.Cells.Item(4,1).Value = "'004237325"
$ws.Cells.Item(4,2).Value = "Ricardo"
$ws.Cells.Item(4,3).Value = 51494.89
$ws.Cells.Item(5,1).Value = "'005277762"
$ws.Cells.Item(5,2).Value = "Nivaldo"
$ws.Cells.Item(5,3).Value = 51452.04
$ws.Cells.Item(6,1).Value = "'004567880"
$ws.Cells.Item(6,2).Value = "Luana"
$ws.Cells.Item(6,3).Value = 32816.84
$ws.Cells.Item(7,1).Value = "'005142611"
$ws.Cells.Item(7,2).Value = "Guilherme"
$ws.Cells.Item(7,3).Value = 22441.24
$ws.Cells.Item(8,1).Value = "'005581299"
$ws.Cells.Item(8,2).Value = "Zilda"
$ws.Cells.Item(8,3).Value = 21847.94
$ws.Cells.Item(9,1).Value = "'004206790"
$ws.Cells.Item(9,2).Value = "Emmanuelle"
$ws.Cells.Item(9,3).Value = 20124.71
$ws.Cells.Item(10,1).Value = "'004481463"
$ws.Cells.Item(10,2).Value = "Mara"
$ws.Cells.Item(10,3).Value = 20000
$ws.Cells.Item(11,1).Value = "'004452476"
$ws.Cells.Item(11,2).Value = "Ivone"
$ws.Cells.Item(11,3).Value = 17502.76
$ws.Cells.Item(12,1).Value = "'004804036"
$ws.Cells.Item(12,2).Value = "Luciana"
$ws.Cells.Item(12,3).Value = 16424.79
$ws.Cells.Item(13,1).Value = "'005142661"
$ws.Cells.Item(13,2).Value = "Sabrina"
$ws.Cells.Item(13,3).Value = 13500.86
$ws.Cells.Item(14,1).Value = "'005046790"
$ws.Cells.Item(14,2).Value = "Beatriz"
$ws.Cells.Item(14,3).Value = 10000
$ws.Cells.Item(15,1).Value = "'008197156"
$ws.Cells.Item(15,2).Value = "Marcio"
$ws.Cells.Item(15,3).Value = 10000
$ws.Cells.Item(16,1).Value = "'005531186"
$ws.Cells.Item(16,2).Value = "Rafael"
$ws.Cells.Item(16,3).Value = 9912
$ws.Cells.Item(17,1).Value = "'004550605"
$ws.Cells.Item(17,2).Value = "Rejane"
$ws.Cells.Item(17,3).Value = 9737.54
$ws.Cells.Item(18,1).Value = "'005070742"
$ws.Cells.Item(18,2).Value = "Juscelino"
$ws.Cells.Item(18,3).Value = 6000
$ws.Cells.Item(19,1).Value = "'008457882"
$ws.Cells.Item(19,2).Value = "Maria"
$ws.Cells.Item(19,3).Value = 5321.99
$ws.Cells.Item(20,1).Value = "'000330949"
$ws.Cells.Item(20,2).Value = "Renato"
$ws.Cells.Item(20,3).Value = 3279.29
$ws.Cells.Item(21,1).Value = "'004504449"
$ws.Cells.Item(21,2).Value = "Kelma"
$ws.Cells.Item(21,3).Value = 2187.85
$ws.Cells.Item(22,1).Value = "'004368468"
$ws.Cells.Item(22,2).Value = "Ahmad"
$ws.Cells.Item(22,3).Value = 1966.84
$ws.Cells.Item(23,1).Value = "'005135281"
$ws.Cells.Item(23,2).Value = "Rafael"
$ws.Cells.Item(23,3).Value = 1500
$ws.Cells.Item(24,1).Value = "'001761119"
$ws.Cells.Item(24,2).Value = "Bluemetrix"
$ws.Cells.Item(24,3).Value = 1053.08
$ws.Cells.Item(25,1).Value = "'004332747"
$ws.Cells.Item(25,2).Value = "Lohran"
$ws.Cells.Item(25,3).Value = 989.79
$ws.Cells.Item(26,1).Value = "'004392159"
$ws.Cells.Item(26,2).Value = "Rodrigo"
$ws.Cells.Item(26,3).Value = 902.31
$ws.Cells.Item(27,1).Value = "'004221454"
$ws.Cells.Item(27,2).Value = "Manuela"
$ws.Cells.Item(27,3).Value = 455.58
$ws.Cells.Item(28,1).Value = "'005637820"
$ws.Cells.Item(28,2).Value = "Guilherme"
$ws.Cells.Item(28,3).Value = 439.12
$ws.Cells.Item(29,1).Value = "'004508526"
$ws.Cells.Item(29,2).Value = "Cassio"
$ws.Cells.Item(29,3).Value = 400.9
$ws.Cells.Item(30,1).Value = "'004458563"
$ws.Cells.Item(30,2).Value = "Luiz"
$ws.Cells.Item(30,3).Value = 278.42
$ws.Cells.Item(31,1).Value = "'004556853"
$ws.Cells.Item(31,2).Value = "Marcel"
$ws.Cells.Item(31,3).Value = 241.93
$ws.Cells.Item(32,1).Value = "'008177213"
$ws.Cells.Item(32,2).Value = "Pedro"
$ws.Cells.Item(32,3).Value = 239.36
$ws.Cells.Item(33,1).Value = "'001719494"
$ws.Cells.Item(33,2).Value = "Luis"
$ws.Cells.Item(33,3).Value = 235.55
$ws.Cells.Item(34,1).Value = "'004261201"
$ws.Cells.Item(34,2).Value = "Ana"
$ws.Cells.Item(34,3).Value = 220
$ws.Cells.Item(35,1).Value = "'004322719"
$ws.Cells.Item(35,2).Value = "Gisela"
$ws.Cells.Item(35,3).Value = 220
$ws.Cells.Item(36,1).Value = "'003512801"
$ws.Cells.Item(36,2).Value = "Lais"
$ws.Cells.Item(36,3).Value = 211.2
$ws.Cells.Item(37,1).Value = "'004359408"
$ws.Cells.Item(37,2).Value = "Hepta"
$ws.Cells.Item(37,3).Value = 187.01
$ws.Cells.Item(38,1).Value = "'004377713"
$ws.Cells.Item(38,2).Value = "Danieli"
$ws.Cells.Item(38,3).Value = 164.55
$ws.Cells.Item(39,1).Value = "'004328934"
$ws.Cells.Item(39,2).Value = "Valeria"
$ws.Cells.Item(39,3).Value = 163.84
$ws.Cells.Item(40,1).Value = "'004405476"
$ws.Cells.Item(40,2).Value = "Mariana"
$ws.Cells.Item(40,3).Value = 98.98
$ws.Cells.Item(41,1).Value = "'004907688"
$ws.Cells.Item(41,2).Value = "Heitor"
$ws.Cells.Item(41,3).Value = 97.41
$ws.Cells.Item(42,1).Value = "'004431591"
$ws.Cells.Item(42,2).Value = "Mario"
$ws.Cells.Item(42,3).Value = 94.24
$ws.Cells.Item(43,1).Value = "'004974089"
$ws.Cells.Item(43,2).Value = "Celia"
$ws.Cells.Item(43,3).Value = 93.23
$ws.Cells.Item(44,1).Value = "'005701765"
$ws.Cells.Item(44,2).Value = "F"
$ws.Cells.Item(44,3).Value = 92.31
$ws.Cells.Item(45,1).Value = "'004340223"
$ws.Cells.Item(45,2).Value = "Pedro"
$ws.Cells.Item(45,3).Value = 88.45
$ws.Cells.Item(46,1).Value = "'004360431"
$ws.Cells.Item(46,2).Value = "Carlos"
$ws.Cells.Item(46,3).Value = 87.52
$ws.Cells.Item(47,1).Value = "'004809902"
$ws.Cells.Item(47,2).Value = "Pedro"
$ws.Cells.Item(47,3).Value = 86.78
$ws.Cells.Item(48,1).Value = "'004212132"
$ws.Cells.Item(48,2).Value = "Joao"
$ws.Cells.Item(48,3).Value = 86.38
$ws.Cells.Item(49,1).Value = "'004384258"
$ws.Cells.Item(49,2).Value = "Paula"
$ws.Cells.Item(49,3).Value = 86.25
$ws.Cells.Item(50,1).Value = "'001368670"
$ws.Cells.Item(50,2).Value = "Thiago"
$ws.Cells.Item(50,3).Value = 85.95
$ws.Cells.Item(51,1).Value = "'000827730"
$ws.Cells.Item(51,2).Value = "Luciana"
$ws.Cells.Item(51,3).Value = 84.62
$ws.Cells.Item(52,1).Value = "'005076418"
$ws.Cells.Item(52,2).Value = "Leonardo"
$ws.Cells.Item(52,3).Value = 81.87
$ws.Cells.Item(53,1).Value = "'004230529"
$ws.Cells.Item(53,2).Value = "Lais"
$ws.Cells.Item(53,3).Value = 72.24
$ws.Cells.Item(54,1).Value = "'004404342"
$ws.Cells.Item(54,2).Value = "Adson"
$ws.Cells.Item(54,3).Value = 54.9
$ws.Cells.Item(55,1).Value = "'004400640"
$ws.Cells.Item(55,2).Value = "Felipe"
$ws.Cells.Item(55,3).Value = 51.44
$ws.Cells.Item(56,1).Value = "'004454365"
$ws.Cells.Item(56,2).Value = "Rafael"
$ws.Cells.Item(56,3).Value = 50.32
$ws.Cells.Item(57,1).Value = "'001731007"
$ws.Cells.Item(57,2).Value = "Guilherme"
$ws.Cells.Item(57,3).Value = 44.79
$ws.Cells.Item(58,1).Value = "'005245032"
$ws.Cells.Item(58,2).Value = "Rosa"
$ws.Cells.Item(58,3).Value = 39.91
$ws.Cells.Item(59,1).Value = "'004238164"
$ws.Cells.Item(59,2).Value = "Daniela"
$ws.Cells.Item(59,3).Value = 39.15
$ws.Cells.Item(60,1).Value = "'008336332"
$ws.Cells.Item(60,2).Value = "Carlos"
$ws.Cells.Item(60,3).Value = 37.15
$ws.Cells.Item(61,1).Value = "'004806286"
$ws.Cells.Item(61,2).Value = "Vera"
$ws.Cells.Item(61,3).Value = 35.8
$ws.Cells.Item(62,1).Value = "'004264780"
$ws.Cells.Item(62,2).Value = "Marcelo"
$ws.Cells.Item(62,3).Value = 35.29
$ws.Cells.Item(63,1).Value = "'004119016"
$ws.Cells.Item(63,2).Value = "Hemat"
$ws.Cells.Item(63,3).Value = 35.24
$ws.Cells.Item(64,1).Value = "'004340036"
$ws.Cells.Item(64,2).Value = "Eduardo"
$ws.Cells.Item(64,3).Value = 33.82
$ws.Cells.Item(65,1).Value = "'004452597"
$ws.Cells.Item(65,2).Value = "Lara"
$ws.Cells.Item(65,3).Value = 31.21
$ws.Cells.Item(66,1).Value = "'002894447"
$ws.Cells.Item(66,2).Value = "Joao"
$ws.Cells.Item(66,3).Value = 31.04
$ws.Cells.Item(67,1).Value = "'004404724"
$ws.Cells.Item(67,2).Value = "Leandro"
$ws.Cells.Item(67,3).Value = 30.02
$ws.Cells.Item(68,1).Value = "'005927101"
$ws.Cells.Item(68,2).Value = "Simone"
$ws.Cells.Item(68,3).Value = 30
$ws.Cells.Item(69,1).Value = "'004377415"
$ws.Cells.Item(69,2).Value = "Angela"
$ws.Cells.Item(69,3).Value = 28.73
$ws.Cells.Item(70,1).Value = "'004389994"
$ws.Cells.Item(70,2).Value = "Polyanna"
$ws.Cells.Item(70,3).Value = 26
$ws.Cells.Item(71,1).Value = "'004350197"
$ws.Cells.Item(71,2).Value = "Gisela"
$ws.Cells.Item(71,3).Value = 25.08
$ws.Cells.Item(72,1).Value = "'005186167"
$ws.Cells.Item(72,2).Value = "Andrea"
$ws.Cells.Item(72,3).Value = 22.15
$ws.Cells.Item(73,1).Value = "'004371857"
$ws.Cells.Item(73,2).Value = "Nazareth"
$ws.Cells.Item(73,3).Value = 21.52
$ws.Cells.Item(74,1).Value = "'004388077"
$ws.Cells.Item(74,2).Value = "Wladmir"
$ws.Cells.Item(74,3).Value = 20.89
$ws.Cells.Item(75,1).Value = "'004214604"
$ws.Cells.Item(75,2).Value = "Maria"
$ws.Cells.Item(75,3).Value = 20.75
$ws.Cells.Item(76,1).Value = "'004204255"
$ws.Cells.Item(76,2).Value = "Amado"
$ws.Cells.Item(76,3).Value = 18.77
$ws.Cells.Item(77,1).Value = "'005374916"
$ws.Cells.Item(77,2).Value = "Marco"
$ws.Cells.Item(77,3).Value = 17.86
$ws.Cells.Item(78,1).Value = "'004422594"
$ws.Cells.Item(78,2).Value = "Wandir"
$ws.Cells.Item(78,3).Value = 14.67
$ws.Cells.Item(79,1).Value = "'005905713"
$ws.Cells.Item(79,2).Value = "Neila"
$ws.Cells.Item(79,3).Value = 13.77
$ws.Cells.Item(80,1).Value = "'005135105"
$ws.Cells.Item(80,2).Value = "Brenner"
$ws.Cells.Item(80,3).Value = 11.26
$ws.Cells.Item(81,1).Value = "'004976625"
$ws.Cells.Item(81,2).Value = "Norton"
$ws.Cells.Item(81,3).Value = 10.76
$ws.Cells.Item(82,1).Value = "'004480134"
$ws.Cells.Item(82,2).Value = "Jose"
$ws.Cells.Item(82,3).Value = 10.61
$ws.Cells.Item(83,1).Value = "'004419141"
$ws.Cells.Item(83,2).Value = "Paulo"
$ws.Cells.Item(83,3).Value = 10.42
$ws.Cells.Item(84,1).Value = "'004420763"
$ws.Cells.Item(84,2).Value = "Christian"
$ws.Cells.Item(84,3).Value = 10.3
$ws.Cells.Item(85,1).Value = "'004216298"
$ws.Cells.Item(85,2).Value = "Flordeliz"
$ws.Cells.Item(85,3).Value = 9.82
$ws.Cells.Item(86,1).Value = "'001294033"
$ws.Cells.Item(86,2).Value = "Viviane"
$ws.Cells.Item(86,3).Value = 8.92
$ws.Cells.Item(87,1).Value = "'008004995"
$ws.Cells.Item(87,2).Value = "Jose"
$ws.Cells.Item(87,3).Value = 8.2
$ws.Cells.Item(88,1).Value = "'005043894"
$ws.Cells.Item(88,2).Value = "Naiara"
$ws.Cells.Item(88,3).Value = 8
$ws.Cells.Item(89,1).Value = "'004530494"
$ws.Cells.Item(89,2).Value = "Rosangela"
$ws.Cells.Item(89,3).Value = 7.05
$ws.Cells.Item(90,1).Value = "'005268516"
$ws.Cells.Item(90,2).Value = "Luis"
$ws.Cells.Item(90,3).Value = 6.35
$ws.Cells.Item(91,1).Value = "'008013889"
$ws.Cells.Item(91,2).Value = "Carolina"
$ws.Cells.Item(91,3).Value = 6.24
$ws.Cells.Item(92,1).Value = "'005905737"
$ws.Cells.Item(92,2).Value = "Cairo"
$ws.Cells.Item(92,3).Value = 6.17
$ws.Cells.Item(93,1).Value = "'004448501"
$ws.Cells.Item(93,2).Value = "Joao"
$ws.Cells.Item(93,3).Value = 5.55
$ws.Cells.Item(94,1).Value = "'004756968"
$ws.Cells.Item(94,2).Value = "Daniely"
$ws.Cells.Item(94,3).Value = 4.66
$ws.Cells.Item(95,1).Value = "'008244362"
$ws.Cells.Item(95,2).Value = "Lincoln"
$ws.Cells.Item(95,3).Value = 3.28
$ws.Cells.Item(96,1).Value = "'004308815"
$ws.Cells.Item(96,2).Value = "Zeli"
$ws.Cells.Item(96,3).Value = 2.51
$ws.Cells.Item(97,1).Value = "'004419765"
$ws.Cells.Item(97,2).Value = "Walter"
$ws.Cells.Item(97,3).Value = 2.41
$ws.Cells.Item(98,1).Value = "'004340217"
$ws.Cells.Item(98,2).Value = "Augusto"
$ws.Cells.Item(98,3).Value = 1.87
$ws.Cells.Item(99,1).Value = "'001882235"
$ws.Cells.Item(99,2).Value = "Lago"
$ws.Cells.Item(99,3).Value = 1.84
$ws.Cells.Item(100,1).Value = "'004840589"
$ws.Cells.Item(100,2).Value = "Leda"
$ws.Cells.Item(100,3).Value = 1.8
$ws.Cells.Item(101,1).Value = "'004460491"
$ws.Cells.Item(101,2).Value = "Pedro"
$ws.Cells.Item(101,3).Value = 1.72
$ws.Cells.Item(102,1).Value = "'005022526"
$ws.Cells.Item(102,2).Value = "Alexandre"
$ws.Cells.Item(102,3).Value = 1.7
$ws.Cells.Item(103,1).Value = "'004220849"
$ws.Cells.Item(103,2).Value = "Dulce"
$ws.Cells.Item(103,3).Value = 1.68
$ws.Cells.Item(104,1).Value = "'004335251"
$ws.Cells.Item(104,2).Value = "Edmundo"
$ws.Cells.Item(104,3).Value = 1.54
$ws.Cells.Item(105,1).Value = "'005575050"
$ws.Cells.Item(105,2).Value = "Garcia"
$ws.Cells.Item(105,3).Value = 1.54
$ws.Cells.Item(106,1).Value = "'004855960"
$ws.Cells.Item(106,2).Value = "Cleria"
$ws.Cells.Item(106,3).Value = 1.51
$ws.Cells.Item(107,1).Value = "'004214460"
$ws.Cells.Item(107,2).Value = "Maria"
$ws.Cells.Item(107,3).Value = 1.5
$ws.Cells.Item(108,1).Value = "'004713953"
$ws.Cells.Item(108,2).Value = "Alessandra"
$ws.Cells.Item(108,3).Value = 1.47
$ws.Cells.Item(109,1).Value = "'004805333"
$ws.Cells.Item(109,2).Value = "Rosana"
$ws.Cells.Item(109,3).Value = 1.47
$ws.Cells.Item(110,1).Value = "'004218542"
$ws.Cells.Item(110,2).Value = "Jose"
$ws.Cells.Item(110,3).Value = 1.45
$ws.Cells.Item(111,1).Value = "'005273382"
$ws.Cells.Item(111,2).Value = "Mvfc"
$ws.Cells.Item(111,3).Value = 1.44
$ws.Cells.Item(112,1).Value = "'004527606"
$ws.Cells.Item(112,2).Value = "Marcia"
$ws.Cells.Item(112,3).Value = 1.43
$ws.Cells.Item(113,1).Value = "'005886225"
$ws.Cells.Item(113,2).Value = "Vinicius"
$ws.Cells.Item(113,3).Value = 1.43
$ws.Cells.Item(114,1).Value = "'004321092"
$ws.Cells.Item(114,2).Value = "Daniel"
$ws.Cells.Item(114,3).Value = 1.37
$ws.Cells.Item(115,1).Value = "'005142624"
$ws.Cells.Item(115,2).Value = "Rodrigo"
$ws.Cells.Item(115,3).Value = 1.36
$ws.Cells.Item(116,1).Value = "'004451652"
$ws.Cells.Item(116,2).Value = "Mateus"
$ws.Cells.Item(116,3).Value = 1.35
$ws.Cells.Item(117,1).Value = "'004587511"
$ws.Cells.Item(117,2).Value = "Carlos"
$ws.Cells.Item(117,3).Value = 1.34
$ws.Cells.Item(118,1).Value = "'004213139"
$ws.Cells.Item(118,2).Value = "Leonardo"
$ws.Cells.Item(118,3).Value = 1.3
$ws.Cells.Item(119,1).Value = "'004381415"
$ws.Cells.Item(119,2).Value = "Joao"
$ws.Cells.Item(119,3).Value = 1.26
$ws.Cells.Item(120,1).Value = "'004482090"
$ws.Cells.Item(120,2).Value = "Cezar"
$ws.Cells.Item(120,3).Value = 1.25
$ws.Cells.Item(121,1).Value = "'005018038"
$ws.Cells.Item(121,2).Value = "Elaine"
$ws.Cells.Item(121,3).Value = 1.25
$ws.Cells.Item(122,1).Value = "'004975924"
$ws.Cells.Item(122,2).Value = "Sergio"
$ws.Cells.Item(122,3).Value = 1.24
$ws.Cells.Item(123,1).Value = "'008328804"
$ws.Cells.Item(123,2).Value = "Sonia"
$ws.Cells.Item(123,3).Value = 1.19
$ws.Cells.Item(124,1).Value = "'004313254"
$ws.Cells.Item(124,2).Value = "Gustavo"
$ws.Cells.Item(124,3).Value = 1.17
$ws.Cells.Item(125,1).Value = "'004911541"
$ws.Cells.Item(125,2).Value = "Tiago"
$ws.Cells.Item(125,3).Value = 1.16
$ws.Cells.Item(126,1).Value = "'004267976"
$ws.Cells.Item(126,2).Value = "E3"
$ws.Cells.Item(126,3).Value = 1.12
$ws.Cells.Item(127,1).Value = "'004479734"
$ws.Cells.Item(127,2).Value = "Rodrigo"
$ws.Cells.Item(127,3).Value = 1.12
$ws.Cells.Item(128,1).Value = "'005381719"
$ws.Cells.Item(128,2).Value = "Maria"
$ws.Cells.Item(128,3).Value = 1.11
$ws.Cells.Item(129,1).Value = "'004693308"
$ws.Cells.Item(129,2).Value = "Laura"
$ws.Cells.Item(129,3).Value = 1.03
$ws.Cells.Item(130,1).Value = "'004462543"
$ws.Cells.Item(130,2).Value = "Rodolfo"
$ws.Cells.Item(130,3).Value = 1.01
$ws.Cells.Item(131,1).Value = "'004360430"
$ws.Cells.Item(131,2).Value = "Viomar"
$ws.Cells.Item(131,3).Value = 1
$ws.Cells.Item(132,1).Value = "'005440756"
$ws.Cells.Item(132,2).Value = "Valeria"
$ws.Cells.Item(132,3).Value = 1
$ws.Cells.Item(133,1).Value = "'004486497"
$ws.Cells.Item(133,2).Value = "Elena"
$ws.Cells.Item(133,3).Value = 0.96
$ws.Cells.Item(134,1).Value = "'004242237"
$ws.Cells.Item(134,2).Value = "Mariah"
$ws.Cells.Item(134,3).Value = 0.95
$ws.Cells.Item(135,1).Value = "'005146441"
$ws.Cells.Item(135,2).Value = "Jose"
$ws.Cells.Item(135,3).Value = 0.95
$ws.Cells.Item(136,1).Value = "'008012870"
$ws.Cells.Item(136,2).Value = "Ana"
$ws.Cells.Item(136,3).Value = 0.92
$ws.Cells.Item(137,1).Value = "'005535788"
$ws.Cells.Item(137,2).Value = "Emilia"
$ws.Cells.Item(137,3).Value = 0.89
$ws.Cells.Item(138,1).Value = "'004115403"
$ws.Cells.Item(138,2).Value = "Hebert"
$ws.Cells.Item(138,3).Value = 0.88
$ws.Cells.Item(139,1).Value = "'001759765"
$ws.Cells.Item(139,2).Value = "Natal"
$ws.Cells.Item(139,3).Value = 0.86
$ws.Cells.Item(140,1).Value = "'005428871"
$ws.Cells.Item(140,2).Value = "Rosangela"
$ws.Cells.Item(140,3).Value = 0.84
$ws.Cells.Item(141,1).Value = "'004497875"
$ws.Cells.Item(141,2).Value = "Henrique"
$ws.Cells.Item(141,3).Value = 0.83
$ws.Cells.Item(142,1).Value = "'005232019"
$ws.Cells.Item(142,2).Value = "Pedro"
$ws.Cells.Item(142,3).Value = 0.83
$ws.Cells.Item(143,1).Value = "'005683532"
$ws.Cells.Item(143,2).Value = "Sylverson"
$ws.Cells.Item(143,3).Value = 0.82
$ws.Cells.Item(144,1).Value = "'004223502"
$ws.Cells.Item(144,2).Value = "Bruna"
$ws.Cells.Item(144,3).Value = 0.78
$ws.Cells.Item(145,1).Value = "'001000288"
$ws.Cells.Item(145,2).Value = "Isabella"
$ws.Cells.Item(145,3).Value = 0.73
$ws.Cells.Item(146,1).Value = "'008298906"
$ws.Cells.Item(146,2).Value = "Fernando"
$ws.Cells.Item(146,3).Value = 0.72
$ws.Cells.Item(147,1).Value = "'004346716"
$ws.Cells.Item(147,2).Value = "Tiago"
$ws.Cells.Item(147,3).Value = 0.71
$ws.Cells.Item(148,1).Value = "'004425261"
$ws.Cells.Item(148,2).Value = "Thaysa"
$ws.Cells.Item(148,3).Value = 0.71
$ws.Cells.Item(149,1).Value = "'004588677"
$ws.Cells.Item(149,2).Value = "Rachel"
$ws.Cells.Item(149,3).Value = 0.71
$ws.Cells.Item(150,1).Value = "'005660155"
$ws.Cells.Item(150,2).Value = "Carolina"
$ws.Cells.Item(150,3).Value = 0.7
$ws.Cells.Item(151,1).Value = "'004335144"
$ws.Cells.Item(151,2).Value = "Edmundo"
$ws.Cells.Item(151,3).Value = 0.67
$ws.Cells.Item(152,1).Value = "'008004835"
$ws.Cells.Item(152,2).Value = "Sergio"
$ws.Cells.Item(152,3).Value = 0.67
$ws.Cells.Item(153,1).Value = "'004473942"
$ws.Cells.Item(153,2).Value = "Daianne"
$ws.Cells.Item(153,3).Value = 0.62
$ws.Cells.Item(154,1).Value = "'005924958"
$ws.Cells.Item(154,2).Value = "Tiago"
$ws.Cells.Item(154,3).Value = 0.6
$ws.Cells.Item(155,1).Value = "'005338054"
$ws.Cells.Item(155,2).Value = "Elaine"
$ws.Cells.Item(155,3).Value = 0.57
$ws.Cells.Item(156,1).Value = "'008119302"
$ws.Cells.Item(156,2).Value = "Vitor"
$ws.Cells.Item(156,3).Value = 0.56
$ws.Cells.Item(157,1).Value = "'002786022"
$ws.Cells.Item(157,2).Value = "Paulo"
$ws.Cells.Item(157,3).Value = 0.55
$ws.Cells.Item(158,1).Value = "'003489079"
$ws.Cells.Item(158,2).Value = "Bruno"
$ws.Cells.Item(158,3).Value = 0.53
$ws.Cells.Item(159,1).Value = "'005105970"
$ws.Cells.Item(159,2).Value = "Vera"
$ws.Cells.Item(159,3).Value = 0.53
$ws.Cells.Item(160,1).Value = "'005133039"
$ws.Cells.Item(160,2).Value = "Paulo"
$ws.Cells.Item(160,3).Value = 0.5
$ws.Cells.Item(161,1).Value = "'004452946"
$ws.Cells.Item(161,2).Value = "Otavio"
$ws.Cells.Item(161,3).Value = 0.49
$ws.Cells.Item(162,1).Value = "'005110894"
$ws.Cells.Item(162,2).Value = "Ana"
$ws.Cells.Item(162,3).Value = 0.49
$ws.Cells.Item(163,1).Value = "'008072033"
$ws.Cells.Item(163,2).Value = "Leticia"
$ws.Cells.Item(163,3).Value = 0.49
$ws.Cells.Item(164,1).Value = "'005558076"
$ws.Cells.Item(164,2).Value = "Alexandre"
$ws.Cells.Item(164,3).Value = 0.47
$ws.Cells.Item(165,1).Value = "'008035153"
$ws.Cells.Item(165,2).Value = "Claudio"
$ws.Cells.Item(165,3).Value = 0.47
$ws.Cells.Item(166,1).Value = "'003641655"
$ws.Cells.Item(166,2).Value = "Marcelo"
$ws.Cells.Item(166,3).Value = 0.45
$ws.Cells.Item(167,1).Value = "'004805133"
$ws.Cells.Item(167,2).Value = "Patricia"
$ws.Cells.Item(167,3).Value = 0.45
$ws.Cells.Item(168,1).Value = "'004767746"
$ws.Cells.Item(168,2).Value = "Isabele"
$ws.Cells.Item(168,3).Value = 0.44
$ws.Cells.Item(169,1).Value = "'005662526"
$ws.Cells.Item(169,2).Value = "Aguinaldo"
$ws.Cells.Item(169,3).Value = 0.44
$ws.Cells.Item(170,1).Value = "'005949170"
$ws.Cells.Item(170,2).Value = "Cintia"
$ws.Cells.Item(170,3).Value = 0.44
$ws.Cells.Item(171,1).Value = "'005141215"
$ws.Cells.Item(171,2).Value = "Karina"
$ws.Cells.Item(171,3).Value = 0.42
$ws.Cells.Item(172,1).Value = "'004432935"
$ws.Cells.Item(172,2).Value = "Jose"
$ws.Cells.Item(172,3).Value = 0.41
$ws.Cells.Item(173,1).Value = "'005216881"
$ws.Cells.Item(173,2).Value = "Renan"
$ws.Cells.Item(173,3).Value = 0.41
$ws.Cells.Item(174,1).Value = "'005341184"
$ws.Cells.Item(174,2).Value = "Breno"
$ws.Cells.Item(174,3).Value = 0.41
$ws.Cells.Item(175,1).Value = "'005530256"
$ws.Cells.Item(175,2).Value = "Carolina"
$ws.Cells.Item(175,3).Value = 0.41
$ws.Cells.Item(176,1).Value = "'004424671"
$ws.Cells.Item(176,2).Value = "Luisa"
$ws.Cells.Item(176,3).Value = 0.4
$ws.Cells.Item(177,1).Value = "'004424761"
$ws.Cells.Item(177,2).Value = "Pedro"
$ws.Cells.Item(177,3).Value = 0.4
$ws.Cells.Item(178,1).Value = "'004563252"
$ws.Cells.Item(178,2).Value = "Fernando"
$ws.Cells.Item(178,3).Value = 0.4
$ws.Cells.Item(179,1).Value = "'005654767"
$ws.Cells.Item(179,2).Value = "Diego"
$ws.Cells.Item(179,3).Value = 0.4
$ws.Cells.Item(180,1).Value = "'004453302"
$ws.Cells.Item(180,2).Value = "Isabella"
$ws.Cells.Item(180,3).Value = 0.39
$ws.Cells.Item(181,1).Value = "'004466350"
$ws.Cells.Item(181,2).Value = "Raquel"
$ws.Cells.Item(181,3).Value = 0.39
$ws.Cells.Item(182,1).Value = "'004397124"
$ws.Cells.Item(182,2).Value = "Murylo"
$ws.Cells.Item(182,3).Value = 0.38
$ws.Cells.Item(183,1).Value = "'004424714"
$ws.Cells.Item(183,2).Value = "Helena"
$ws.Cells.Item(183,3).Value = 0.38
$ws.Cells.Item(184,1).Value = "'002973105"
$ws.Cells.Item(184,2).Value = "Darlan"
$ws.Cells.Item(184,3).Value = 0.37
$ws.Cells.Item(185,1).Value = "'004972070"
$ws.Cells.Item(185,2).Value = "Maria"
$ws.Cells.Item(185,3).Value = 0.37
$ws.Cells.Item(186,1).Value = "'005725431"
$ws.Cells.Item(186,2).Value = "Bruno"
$ws.Cells.Item(186,3).Value = 0.37
$ws.Cells.Item(187,1).Value = "'003836362"
$ws.Cells.Item(187,2).Value = "Isabella"
$ws.Cells.Item(187,3).Value = 0.33
$ws.Cells.Item(188,1).Value = "'004556150"
$ws.Cells.Item(188,2).Value = "Marina"
$ws.Cells.Item(188,3).Value = 0.33
$ws.Cells.Item(189,1).Value = "'004890544"
$ws.Cells.Item(189,2).Value = "Assako"
$ws.Cells.Item(189,3).Value = 0.3
$ws.Cells.Item(190,1).Value = "'005379541"
$ws.Cells.Item(190,2).Value = "Jose"
$ws.Cells.Item(190,3).Value = 0.28
$ws.Cells.Item(191,1).Value = "'004908680"
$ws.Cells.Item(191,2).Value = "Elene"
$ws.Cells.Item(191,3).Value = 0.26
$ws.Cells.Item(192,1).Value = "'005270025"
$ws.Cells.Item(192,2).Value = "Denize"
$ws.Cells.Item(192,3).Value = 0.25
$ws.Cells.Item(193,1).Value = "'008008723"
$ws.Cells.Item(193,2).Value = "Redrau"
$ws.Cells.Item(193,3).Value = 0.25
$ws.Cells.Item(194,1).Value = "'004207278"
$ws.Cells.Item(194,2).Value = "Cesar"
$ws.Cells.Item(194,3).Value = 0.23
$ws.Cells.Item(195,1).Value = "'005035754"
$ws.Cells.Item(195,2).Value = "Jose"
$ws.Cells.Item(195,3).Value = 0.23
$ws.Cells.Item(196,1).Value = "'005092207"
$ws.Cells.Item(196,2).Value = "Bruno"
$ws.Cells.Item(196,3).Value = 0.23
$ws.Cells.Item(197,1).Value = "'005110867"
$ws.Cells.Item(197,2).Value = "Dig"
$ws.Cells.Item(197,3).Value = 0.23
$ws.Cells.Item(198,1).Value = "'004278033"
$ws.Cells.Item(198,2).Value = "Daisy"
$ws.Cells.Item(198,3).Value = 0.21
$ws.Cells.Item(199,1).Value = "'004612043"
$ws.Cells.Item(199,2).Value = "Yuri"
$ws.Cells.Item(199,3).Value = 0.21
$ws.Cells.Item(200,1).Value = "'004848843"
$ws.Cells.Item(200,2).Value = "Alaercio"
$ws.Cells.Item(200,3).Value = 0.21
$ws.Cells.Item(201,1).Value = "'005697554"
$ws.Cells.Item(201,2).Value = "Maria"
$ws.Cells.Item(201,3).Value = 0.2
$ws.Cells.Item(202,1).Value = "'005304669"
$ws.Cells.Item(202,2).Value = "Emilson"
$ws.Cells.Item(202,3).Value = 0.18
$ws.Cells.Item(203,1).Value = "'004432455"
$ws.Cells.Item(203,2).Value = "Luciana"
$ws.Cells.Item(203,3).Value = 0.17
$ws.Cells.Item(204,1).Value = "'002687737"
$ws.Cells.Item(204,2).Value = "Jose"
$ws.Cells.Item(204,3).Value = 0.16
$ws.Cells.Item(205,1).Value = "'002694089"
$ws.Cells.Item(205,2).Value = "Vitor"
$ws.Cells.Item(205,3).Value = 0.16
$ws.Cells.Item(206,1).Value = "'004339183"
$ws.Cells.Item(206,2).Value = "Jalison"
$ws.Cells.Item(206,3).Value = 0.15
$ws.Cells.Item(207,1).Value = "'004357159"
$ws.Cells.Item(207,2).Value = "Joao"
$ws.Cells.Item(207,3).Value = 0.15
$ws.Cells.Item(208,1).Value = "'004398253"
$ws.Cells.Item(208,2).Value = "Euler"
$ws.Cells.Item(208,3).Value = 0.15
$ws.Cells.Item(209,1).Value = "'004754056"
$ws.Cells.Item(209,2).Value = "Bruno"
$ws.Cells.Item(209,3).Value = 0.15
$ws.Cells.Item(210,1).Value = "'004805269"
$ws.Cells.Item(210,2).Value = "Clisia"
$ws.Cells.Item(210,3).Value = 0.15
$ws.Cells.Item(211,1).Value = "'004320840"
$ws.Cells.Item(211,2).Value = "Natalia"
$ws.Cells.Item(211,3).Value = 0.14
$ws.Cells.Item(212,1).Value = "'005075382"
$ws.Cells.Item(212,2).Value = "Nayara"
$ws.Cells.Item(212,3).Value = 0.14
$ws.Cells.Item(213,1).Value = "'008336728"
$ws.Cells.Item(213,2).Value = "Arthur"
$ws.Cells.Item(213,3).Value = 0.14
$ws.Cells.Item(214,1).Value = "'000938440"
$ws.Cells.Item(214,2).Value = "Base"
$ws.Cells.Item(214,3).Value = 0.12
$ws.Cells.Item(215,1).Value = "'004382374"
$ws.Cells.Item(215,2).Value = "Theomar"
$ws.Cells.Item(215,3).Value = 0.12
$ws.Cells.Item(216,1).Value = "'004551472"
$ws.Cells.Item(216,2).Value = "Diego"
$ws.Cells.Item(216,3).Value = 0.12
$ws.Cells.Item(217,1).Value = "'005313179"
$ws.Cells.Item(217,2).Value = "Maira"
$ws.Cells.Item(217,3).Value = 0.12
$ws.Cells.Item(218,1).Value = "'008071998"
$ws.Cells.Item(218,2).Value = "Isadora"
$ws.Cells.Item(218,3).Value = 0.11
$ws.Cells.Item(219,1).Value = "'004281300"
$ws.Cells.Item(219,2).Value = "Franklin"
$ws.Cells.Item(219,3).Value = 0.1
$ws.Cells.Item(220,1).Value = "'005880251"
$ws.Cells.Item(220,2).Value = "Luiz"
$ws.Cells.Item(220,3).Value = 0.1
$ws.Cells.Item(221,1).Value = "'008351535"
$ws.Cells.Item(221,2).Value = "David"
$ws.Cells.Item(221,3).Value = 0.1
$ws.Cells.Item(222,1).Value = "'004451996"
$ws.Cells.Item(222,2).Value = "Adriano"
$ws.Cells.Item(222,3).Value = 0.09
$ws.Cells.Item(223,1).Value = "'004493324"
$ws.Cells.Item(223,2).Value = "Daniel"
$ws.Cells.Item(223,3).Value = 0.09
$ws.Cells.Item(224,1).Value = "'005047946"
$ws.Cells.Item(224,2).Value = "Gabriel"
$ws.Cells.Item(224,3).Value = 0.09
$ws.Cells.Item(225,1).Value = "'008407512"
$ws.Cells.Item(225,2).Value = "Mauricio"
$ws.Cells.Item(225,3).Value = 0.09
$ws.Cells.Item(226,1).Value = "'004213373"
$ws.Cells.Item(226,2).Value = "Alexandre"
$ws.Cells.Item(226,3).Value = 0.08
$ws.Cells.Item(227,1).Value = "'004216434"
$ws.Cells.Item(227,2).Value = "Jaime"
$ws.Cells.Item(227,3).Value = 0.08
$ws.Cells.Item(228,1).Value = "'005720146"
$ws.Cells.Item(228,2).Value = "Jordanna"
$ws.Cells.Item(228,3).Value = 0.07
$ws.Cells.Item(229,1).Value = "'004691225"
$ws.Cells.Item(229,2).Value = "Anna"
$ws.Cells.Item(229,3).Value = 0.06
$ws.Cells.Item(230,1).Value = "'005171652"
$ws.Cells.Item(230,2).Value = "Bruno"
$ws.Cells.Item(230,3).Value = 0.06
$ws.Cells.Item(231,1).Value = "'008026930"
$ws.Cells.Item(231,2).Value = "Joao"
$ws.Cells.Item(231,3).Value = 0.06
$ws.Cells.Item(232,1).Value = "'008037529"
$ws.Cells.Item(232,2).Value = "Melissa"
$ws.Cells.Item(232,3).Value = 0.05
$ws.Cells.Item(233,1).Value = "'008123677"
$ws.Cells.Item(233,2).Value = "Priscilla"
$ws.Cells.Item(233,3).Value = 0.05
$ws.Cells.Item(234,1).Value = "'005274028"
$ws.Cells.Item(234,2).Value = "Rafael"
$ws.Cells.Item(234,3).Value = 0.04
$ws.Cells.Item(235,1).Value = "'003107135"
$ws.Cells.Item(235,2).Value = "Ana"
$ws.Cells.Item(235,3).Value = 0.03
$ws.Cells.Item(236,1).Value = "'004329229"
$ws.Cells.Item(236,2).Value = "Gabriel"
$ws.Cells.Item(236,3).Value = 0.03
$ws.Cells.Item(237,1).Value = "'004946997"
$ws.Cells.Item(237,2).Value = "Eduardo"
$ws.Cells.Item(237,3).Value = 0.03
$ws.Cells.Item(238,1).Value = "'008070544"
$ws.Cells.Item(238,2).Value = "Marina"
$ws.Cells.Item(238,3).Value = 0.03
$ws.Cells.Item(239,1).Value = "'008110684"
$ws.Cells.Item(239,2).Value = "Edval"
$ws.Cells.Item(239,3).Value = 0.03
$ws.Cells.Item(240,1).Value = "'008365179"
$ws.Cells.Item(240,2).Value = "Giovana"
$ws.Cells.Item(240,3).Value = 0.03
$ws.Cells.Item(241,1).Value = "'004384131"
$ws.Cells.Item(241,2).Value = "Andre"
$ws.Cells.Item(241,3).Value = 0.02
$ws.Cells.Item(242,1).Value = "'004457389"
$ws.Cells.Item(242,2).Value = "Rafael"
$ws.Cells.Item(242,3).Value = 0.02
$ws.Cells.Item(243,1).Value = "'005295509"
$ws.Cells.Item(243,2).Value = "Bhruna"
$ws.Cells.Item(243,3).Value = 0.02
$ws.Cells.Item(244,1).Value = "'008090243"
$ws.Cells.Item(244,2).Value = "Gabriel"
$ws.Cells.Item(244,3).Value = 0.02
$ws.Cells.Item(245,1).Value = "'002878817"
$ws.Cells.Item(245,2).Value = "Guilherme"
$ws.Cells.Item(245,3).Value = 0.01
$ws.Cells.Item(246,1).Value = "'004223226"
$ws.Cells.Item(246,2).Value = "Yeshua"
$ws.Cells.Item(246,3).Value = 0.01
$ws.Cells.Item(247,1).Value = "'004272426"
$ws.Cells.Item(247,2).Value = "Rodrigo"
$ws.Cells.Item(247,3).Value = 0.01
$ws.Cells.Item(248,1).Value = "'004400000"
$ws.Cells.Item(248,2).Value = "Vilma"
$ws.Cells.Item(248,3).Value = 0.01
$ws.Cells.Item(249,1).Value = "'004462930"
$ws.Cells.Item(249,2).Value = "Walter"
$ws.Cells.Item(249,3).Value = 0.01
$ws.Cells.Item(250,1).Value = "'004472076"
$ws.Cells.Item(250,2).Value = "Rubens"
$ws.Cells.Item(250,3).Value = 0.01
$ws.Cells.Item(251,1).Value = "'004862746"
$ws.Cells.Item(251,2).Value = "Cesar"
$ws.Cells.Item(251,3).Value = 0.01
$ws.Cells.Item(252,1).Value = "'005068961"
$ws.Cells.Item(252,2).Value = "Jorge"
$ws.Cells.Item(252,3).Value = 0.01
$ws.Cells.Item(253,1).Value = "'005142592"
$ws.Cells.Item(253,2).Value = "Alberto"
$ws.Cells.Item(253,3).Value = 0.01

# --- Blank separator row (clear all three columns; this row previously held data) ---
$ws.Cells.Item(254,1).Value = ""
$ws.Cells.Item(254,2).Value = ""
$ws.Cells.Item(254,3).Value = ""

# --- Footer / filter-notes row (clear columns B & C; this row previously held data) ---
$ws.Cells.Item(255,1).Value = "Filtros aplicados:`nnr_saldo_disponivel não é 0`nPosição é Posição D-1`nCARTEIRA não está em branco`ntela é Financeiro`ntipo_conta é Outros`nNR_CONTA não está em branco`nTIPO_LANCAMENTO não é ED, ET ou Liquidação Doador`nCARTEIRA não é (Em branco)"
$ws.Cells.Item(255,2).Value = ""
$ws.Cells.Item(255,3).Value = ""

# --- The new layout is shorter than the old one (net removal of 3 data rows),
#     so clear out the now-stale tail rows (256-258) left over from before. ---
if (258 -gt 255) {
    $ws.Range("A256:C258").EntireRow.Delete()
}
